# Applies the change described by the diff:
#  1. Adds a new worksheet "ODI Batting Extra" (sheetId=4) after "ODI Bowling".
#  2. Populates it with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns, copying the header style
#     used on the other data sheets.
#  3. Clears the (already-empty) INNING_NUMBER cells B2, B4, B5, B6, B7 on the
#     "ODI Batting" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clear empty inning-number cells on "ODI Batting" (B2, B4, B5, B6, B7)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
foreach ($r in @(2, 4, 5, 6, 7)) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" sheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row values
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold/border/centered header style used by the other sheets
$battingSheet.Range("A1:D1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# MATCH_CODE column (A2:A9) matches MATCH_CODE (D2:D9) already on "ODI Batting"
# -- copy it wholesale so values come across as plain text with the default
# (unstyled) cell format, exactly like the source column.
$battingSheet.Range("D2:D9").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4104)  # xlPasteAll

# Helper: write a numeric-looking (or percent-looking) string as genuine TEXT
# without leaving a lingering explicit cell style behind. Excel auto-detects
# digit/percent strings as numbers unless the cell is pre-formatted as text;
# clearing the format again afterwards drops the cell back to the default
# (unstyled) xf, matching the source workbook's plain inlineStr cells.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Data rows: MATCH_CODE(already filled above), BATTING_POSITION, NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @{ R = 2; B = 9;     C = $null; D = $null; E = $null;   F = "NO" },
    @{ R = 3; B = 10;    C = "0";   D = "0";   E = "1.22%"; F = "NO" },
    @{ R = 4; B = $null; C = $null; D = $null; E = $null;   F = "NO" },
    @{ R = 5; B = $null; C = $null; D = $null; E = $null;   F = "NO" },
    @{ R = 6; B = $null; C = $null; D = $null; E = $null;   F = "NO" },
    @{ R = 7; B = $null; C = $null; D = $null; E = $null;   F = "NO" },
    @{ R = 8; B = 9;     C = "0";   D = "0";   E = "0.34%"; F = "NO" },
    @{ R = 9; B = $null; C = $null; D = $null; E = $null;   F = $null }
)

foreach ($row in $rows) {
    $r = $row.R

    # BATTING_POSITION: real number when present, otherwise blank (text-typed
    # empty string, represented as a formula since COM collapses a literal
    # empty-string assignment into "no cell").
    if ($null -ne $row.B) {
        $newSheet.Cells.Item($r, 2).Value = $row.B
    } else {
        $newSheet.Cells.Item($r, 2).Formula = '=""'
    }

    # NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL: text-typed numeric-looking
    # strings, or blank.
    if ($null -ne $row.C) { Set-TextValue $newSheet.Cells.Item($r, 3) $row.C } else { $newSheet.Cells.Item($r, 3).Formula = '=""' }
    if ($null -ne $row.D) { Set-TextValue $newSheet.Cells.Item($r, 4) $row.D } else { $newSheet.Cells.Item($r, 4).Formula = '=""' }
    if ($null -ne $row.E) { Set-TextValue $newSheet.Cells.Item($r, 5) $row.E } else { $newSheet.Cells.Item($r, 5).Formula = '=""' }

    # MAN_OF_MATCH: plain text ("NO"), or blank for the last row.
    if ($null -ne $row.F) { $newSheet.Cells.Item($r, 6).Value = $row.F } else { $newSheet.Cells.Item($r, 6).Formula = '=""' }
}
